{"js": "// Replace the date line and every \"A\u00d7B=C\" equation cell with its updated value.\n// All old values are unique in the document, so a targeted search-and-replace\n// (one hit expected per lookup) for each old->new pair reproduces the diff\n// exactly, regardless of which table row/cell each value lives in.\nconst replacements = [\n  [\"2025-05-29 Thursday\", \"2025-05-30 Friday\"],\n  [\"11\u00d743=473\", \"33\u00d776=2508\"],\n  [\"43\u00d725=1075\", \"69\u00d725=1725\"],\n  [\"75\u00d791=6825\", \"96\u00d721=2016\"],\n  [\"43\u00d716=688\", \"94\u00d764=6016\"],\n  [\"72\u00d712=864\", \"36\u00d767=2412\"],\n  [\"50\u00d722=1100\", \"43\u00d780=3440\"],\n  [\"75\u00d779=5925\", \"42\u00d737=1554\"],\n  [\"87\u00d742=3654\", \"41\u00d728=1148\"],\n  [\"97\u00d755=5335\", \"28\u00d734=952\"],\n  [\"11\u00d750=550\", \"25\u00d758=1450\"],\n  [\"53\u00d775=3975\", \"28\u00d734=952\"],\n  [\"51\u00d772=3672\", \"33\u00d729=957\"],\n  [\"37\u00d721=777\", \"67\u00d765=4355\"],\n  [\"47\u00d773=3431\", \"18\u00d778=1404\"],\n  [\"78\u00d736=2808\", \"70\u00d779=5530\"],\n  [\"81\u00d773=5913\", \"11\u00d769=759\"],\n  [\"48\u00d740=1920\", \"29\u00d742=1218\"],\n  [\"75\u00d760=4500\", \"65\u00d779=5135\"],\n  [\"77\u00d779=6083\", \"23\u00d776=1748\"],\n  [\"74\u00d730=2220\", \"61\u00d789=5429\"],\n  [\"91\u00d765=5915\", \"79\u00d774=5846\"],\n  [\"68\u00d747=3196\", \"18\u00d738=684\"],\n  [\"69\u00d795=6555\", \"48\u00d732=1536\"],\n  [\"94\u00d785=7990\", \"44\u00d767=2948\"],\n  [\"96\u00d765=6240\", \"27\u00d714=378\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace only the first occurrence (old values are unique in this document).\n  found.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and every \"A\u00d7B=C\" equation cell with its updated\n# value. All old values are unique in the document, so a targeted\n# Find/Replace (one hit expected per lookup, wdReplaceOne) for each\n# old->new pair reproduces the diff exactly, regardless of which table\n# row/cell each value lives in.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-29 Thursday\", \"2025-05-30 Friday\"),\n    @(\"11\u00d743=473\", \"33\u00d776=2508\"),\n    @(\"43\u00d725=1075\", \"69\u00d725=1725\"),\n    @(\"75\u00d791=6825\", \"96\u00d721=2016\"),\n    @(\"43\u00d716=688\", \"94\u00d764=6016\"),\n    @(\"72\u00d712=864\", \"36\u00d767=2412\"),\n    @(\"50\u00d722=1100\", \"43\u00d780=3440\"),\n    @(\"75\u00d779=5925\", \"42\u00d737=1554\"),\n    @(\"87\u00d742=3654\", \"41\u00d728=1148\"),\n    @(\"97\u00d755=5335\", \"28\u00d734=952\"),\n    @(\"11\u00d750=550\", \"25\u00d758=1450\"),\n    @(\"53\u00d775=3975\", \"28\u00d734=952\"),\n    @(\"51\u00d772=3672\", \"33\u00d729=957\"),\n    @(\"37\u00d721=777\", \"67\u00d765=4355\"),\n    @(\"47\u00d773=3431\", \"18\u00d778=1404\"),\n    @(\"78\u00d736=2808\", \"70\u00d779=5530\"),\n    @(\"81\u00d773=5913\", \"11\u00d769=759\"),\n    @(\"48\u00d740=1920\", \"29\u00d742=1218\"),\n    @(\"75\u00d760=4500\", \"65\u00d779=5135\"),\n    @(\"77\u00d779=6083\", \"23\u00d776=1748\"),\n    @(\"74\u00d730=2220\", \"61\u00d789=5429\"),\n    @(\"91\u00d765=5915\", \"79\u00d774=5846\"),\n    @(\"68\u00d747=3196\", \"18\u00d738=684\"),\n    @(\"69\u00d795=6555\", \"48\u00d732=1536\"),\n    @(\"94\u00d785=7990\", \"44\u00d767=2948\"),\n    @(\"96\u00d765=6240\", \"27\u00d714=378\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # wdFindContinue = 1, wdReplaceOne = 1 (replace only the single, unique hit)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1) | Out-Null\n}\n"}
